$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New column J: "Pct >=25 w HS Diploma"
# ---------------------------------------------------------------------------

# Header cell - reuse the same header look as the rest of row 1 (right/center
# aligned, default font/number format) by copying the format already used on
# the adjacent header cell.
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("J1").Value = "Pct >=25 w HS Diploma"

# Register Excel's built-in "Comma" cell style (adds the numFmt 43 / font /
# cellStyleXfs / cellStyles bookkeeping that a real "Apply Cell Styles >
# Comma" action produces) using a scratch row far outside the used range,
# then discard that scratch row completely so no trace of it remains.
$scratch = $ws.Range("A1048576")
$scratch.Style = "Comma"
$scratch.EntireRow.Delete()

# Give the new data cells the same visual result as Excel's "Comma" style
# (thousands separator, 2 decimals, parentheses for negatives) left-aligned
# and vertically centered, matching the rest of the data rows.
$data = $ws.Range("J2:J4")
$data.NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"
$data.HorizontalAlignment = -4131  # xlLeft
$data.VerticalAlignment = -4108    # xlCenter

$ws.Range("J2").Value = 78.431446683554512
$ws.Range("J3").Value = 84.104320036523433
$ws.Range("J4").Value = 87.139870353327197

# Column J width, sized to fit the new header text.
$ws.Range("J1").EntireColumn.ColumnWidth = 20.09765625
